$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 227; this shifts rows 227:322 down to 228:323
# (including copying the "D" column date-number-format style to the new row),
# and the sheet's used-range dimension extends to R323 automatically.
$ws.Rows("227:227").Insert()

# Populate the newly inserted row 227 with its new data.
$ws.Range("A227").Value = 10
$ws.Range("B227").Value = "Vega Modelo de Temuco"
$ws.Range("C227").Value = "La Araucanía"
$ws.Range("D227").Value = 44489
$ws.Range("E227").Value = 9
$ws.Range("F227").Value = 100112043
$ws.Range("G227").Value = "Pepino ensalada"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 110
$ws.Range("K227").Value = 12000
$ws.Range("L227").Value = 12000
$ws.Range("M227").Value = 12000
$ws.Range("N227").Value = "`$/caja 60 unidades"
$ws.Range("O227").Value = "Región de Arica y Parinacota"
$ws.Range("P227").Value = 200
$ws.Range("Q227").Value = 60
$ws.Range("R227").Value = "Hortaliza"
